$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1398
$ws1.Range("F5").Value = 113
$ws1.Range("F7").Value = 11860
$ws1.Range("F8").Value = 4431
$ws1.Range("F13").Value = 2560
$ws1.Range("F16").Value = 50
$ws1.Range("F17").Value = 5144
$ws1.Range("F19").Value = 191
$ws1.Range("F20").Value = 530
$ws1.Range("F21").Value = 11381
$ws1.Range("F22").Value = 11354
$ws1.Range("F23").Value = 20
$ws1.Range("F25").Value = 14

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 3

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1398
$ws4.Range("F5").Value = 113
$ws4.Range("F7").Value = 11860
$ws4.Range("F8").Value = 4431
$ws4.Range("F13").Value = 2560
$ws4.Range("F14").Value = 3
$ws4.Range("F17").Value = 50
$ws4.Range("F18").Value = 5144
$ws4.Range("F20").Value = 191
$ws4.Range("F21").Value = 530
$ws4.Range("F22").Value = 11381
$ws4.Range("F23").Value = 11354
$ws4.Range("F24").Value = 20
$ws4.Range("F26").Value = 14
